# Performance Data.xlsx - add a second "Intel Core i5 @ 2.6 GHz, 8 GB RAM"
# benchmark block (rows 18-28), mirroring the existing block in rows 3-13,
# and move the current selection/scroll position down to the new block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Replicate formatting from the first benchmark block down to the new one ---

# Header band (merged cell with machine name) F3:H3 -> F18:H18
$ws.Range("F3:H3").Copy()
$ws.Range("F18:H18").PasteSpecial(-4122)   # xlPasteFormats

# Column headers ("Time (s)" / "Packet/s") F4:G4 -> F19:G19
$ws.Range("F4:G4").Copy()
$ws.Range("F19:G19").PasteSpecial(-4122)

# Data rows F5:G13 -> F20:G28
$ws.Range("F5:G13").Copy()
$ws.Range("F20:G28").PasteSpecial(-4122)

# Side annotation cells (column H) - copy format individually so that rows
# which have no H cell in the source stay untouched in the destination.
$ws.Range("H6").Copy()
$ws.Range("H21").PasteSpecial(-4122)

$ws.Range("H9").Copy()
$ws.Range("H24").PasteSpecial(-4122)

$ws.Range("H12").Copy()
$ws.Range("H27").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# The old placeholder (blank, style-only) cells in column H for rows that
# must end up with no H cell at all once the new block is filled in.
$ws.Range("H20").Clear()
$ws.Range("H23").Clear()
$ws.Range("H26").Clear()

# --- 2. Merge the new header band, like the source F3:H3 ---
$ws.Range("F18:H18").Merge()

# --- 3. Header text for the new block ---
$ws.Range("F18").Value = "Intel Core i5 @ 2.6 GHz, 8 GB RAM"
$ws.Range("F19").Value = "Time (s)"
$ws.Range("G19").Value = "Packet/s"

# --- 4. Data values (Count column) for the new block ---
$ws.Range("F20").Value = 49
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 41
$ws.Range("F23").Value = 2
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = 22
$ws.Range("F27").Value = 20
$ws.Range("F28").Value = 25

# --- 5. Formulas (Packet/s column) - single assignment over the whole range
#     creates one shared formula group, same as the original block. ---
$ws.Range("G20").Formula = "=C5/F20"
$ws.Range("G21:G28").Formula = "=C6/F21"

# --- 6. Side annotation text / formulas (column H) ---
$ws.Range("H21").Value = "Packet/s AVG time (s)"
$ws.Range("H24").Value = "Process time (s)"
$ws.Range("H27").Value = "Packet/s AVG time (s)"

$ws.Range("H22").Formula = "=AVERAGEA(G20:G28)"
$ws.Range("H25").Formula = "=SUM(F20:F28)"
$ws.Range("H28").Formula = "=D16/H25"

# Formulas with a division operator pick up an automatic number format in
# this engine (unlike the source file, saved by real Excel) - strip that
# back off so the cells stay format-less, like H7/H10/H13 in the source.
$ws.Range("H25").ClearFormats()
$ws.Range("H28").ClearFormats()

# --- 7. View state: scroll so row 7 is the top-visible row, and move the
#     active selection to the cell next to the new block. ---
$ws.Range("I25").Select()
$excel.ActiveWindow.ScrollRow = 7
